# lab_3: "Fix lab_3, add task_6"
#
# This script reproduces (to the extent the COM surface allows) the changes
# described by the authoritative diff:
#   1. A new (empty-looking, single-space) shared string is introduced and
#      used in a brand-new row 19 (cell H19) - this is the "add task_6" bit.
#   2. The sign of the denominator in the O28:O37 formulas flips
#      ((-0.4+1.21*N)/0.19  ->  (-0.4+1.21*N)/-0.19), which flips the sign
#      of every cached/calculated value in that column (and, in real Excel,
#      the chart series that plots it).
#   3. Three stray/duplicate cell formats (on F37, F47, F57) get tidied up
#      to match the formatting already used by their row neighbours.
#   4. One chart marker shrinks from size 7 to size 4.
#   5. The saved selection moves to N23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. New row 19, single space value -> introduces the 14th shared string.
# ---------------------------------------------------------------------
$ws.Range("H19").Value = " "

# ---------------------------------------------------------------------
# 2. Flip the sign of the divisor for the whole O28:O37 shared formula.
#    Writing the anchor cell first and then the rest of the range keeps
#    the shared-formula grouping (O29:O37) intact, same as the original.
# ---------------------------------------------------------------------
$ws.Range("O28").Formula = "=(-0.4+1.21*N28)/-0.19"
$ws.Range("O29:O37").Formula = "=(-0.4+1.21*N29)/-0.19"

# ---------------------------------------------------------------------
# 3. Clean up the formatting on F37 / F47 / F57 so each matches the style
#    already used by the rest of its row (G37/G47/G57).
# ---------------------------------------------------------------------
$ws.Range("G37").Copy()
$ws.Range("F37").PasteSpecial(-4122)

$ws.Range("G47").Copy()
$ws.Range("F47").PasteSpecial(-4122)

$ws.Range("G57").Copy()
$ws.Range("F57").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Shrink the third series' marker on the first chart (7 -> 4).
# ---------------------------------------------------------------------
$chart1 = $ws.ChartObjects(1).Chart
$chart1.SeriesCollection(3).MarkerSize = 4

# ---------------------------------------------------------------------
# 5. Move the live selection to N23 (also scrolls the view near it).
# ---------------------------------------------------------------------
$null = $ws.Range("N23").Select()
